# Updates the "cryptos" worksheet with refreshed price/volume data
# (GitHub Actions scheduled data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.521.79"
$ws.Range("E2").Value = "  -3.84%  "
$ws.Range("D3").Value = "2.536.46"
$ws.Range("E3").Value = "  -3.71%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.46%  "
$ws.Range("D9").Value = "2.542.40"
$ws.Range("E9").Value = "  -3.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.67%  "
$ws.Range("E11").Value = "  -6.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.331"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.40%  "
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "2.983.19"
$ws.Range("E14").Value = "  -3.66%  "
$ws.Range("D15").Value = "58.524.01"
$ws.Range("E15").Value = "  -3.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.26%  "
$ws.Range("E17").Value = "  -6.00%  "
$ws.Range("D18").Value = "2.541.23"
$ws.Range("E18").Value = "  -3.94%  "
$ws.Range("E19").Value = "  -4.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.03%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.410"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -4.84%  "
$ws.Range("D28").Value = "2.652.77"
$ws.Range("E28").Value = "  -3.58%  "
$ws.Range("D29").Value = "0.0₃0785"
$ws.Range("E29").Value = "  -9.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.26%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.913"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.08%  "
$ws.Range("E38").Value = "  -7.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.821"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.47%  "
$ws.Range("E41").Value = "  -7.00%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.83%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "282.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0994"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.77%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.599"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0532"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0226"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.98%  "
$ws.Range("E51").Value = "  -8.25%  "
